$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: extend the width=9 custom sizing from column C to columns C:E ---
$ws.Columns("D:E").ColumnWidth = 8.17

# --- Row 1: weekday date headers across C1:L1 ---
$ws.Range("C1").Value = 42305
$ws.Range("D1").Value = 42306
$ws.Range("E1").Value = 42307
$ws.Range("F1").Value = 42310
$ws.Range("G1").Value = 42311
$ws.Range("H1").Value = 42312
$ws.Range("I1").Value = 42313
$ws.Range("J1").Value = 42314
$ws.Range("K1").Value = 42317
$ws.Range("L1").Value = 42318

$ws.Range("C1:L1").NumberFormat = "d-mmm"
$ws.Range("C1:L1").HorizontalAlignment = -4108

# H1 picks up the same font used by F4 (distinct from the rest of the row)
$ws.Range("F4").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").NumberFormat = "d-mmm"
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").VerticalAlignment = -4107
$ws.Range("H1").Value = 42312

# --- Row 10 updates ---
$ws.Range("F10").Value = 0.5

$ws.Range("F4").Copy()
$ws.Range("G10").PasteSpecial(-4122)
$ws.Range("G10").Value = 0

$ws.Range("C10").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$ws.Range("I10").Value = 1.5

# --- Selection ---
$ws.Range("I17").Select()
